# Auto-generated: apply Leve profit data refresh per diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 50.6
$ws.Range("I5").Value = 46.875
$ws.Range("J5").Value = 65.5
$ws.Range("K5").Value = 46.875
$ws.Range("L5").Value = 65.5
$ws.Range("M5").Value = 68.125
$ws.Range("N5").Value = -295.5
$ws.Range("H12").Value = 242.33333
$ws.Range("I12").Value = 242.33333
$ws.Range("K12").Value = 242.33333
$ws.Range("M12").Value = -72.33332999999999
$ws.Range("H111").Value = 1041.1
$ws.Range("I111").Value = 1045.6666
$ws.Range("K111").Value = 3136.9998
$ws.Range("M111").Value = -69.99980000000005
$ws.Range("H132").Value = 6747.0435
$ws.Range("I132").Value = 5394.45
$ws.Range("J132").Value = 15764.333
$ws.Range("K132").Value = 16183.35
$ws.Range("L132").Value = 47292.999
$ws.Range("M132").Value = -13653.35
$ws.Range("N132").Value = -52352.999
$ws.Range("H138").Value = 4677.4
$ws.Range("J138").Value = 5908.1924
$ws.Range("L138").Value = 17724.5772
$ws.Range("N138").Value = -28004.5772
$ws.Range("H141").Value = 4534.5454
$ws.Range("I141").Value = 4534.5454
$ws.Range("K141").Value = 13603.6362
$ws.Range("M141").Value = -8423.636200000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 25029
$ws.Range("J36").Value = 25029
$ws.Range("L36").Value = 25029
$ws.Range("N36").Value = -25721
$ws.Range("H97").Value = 5822.5356
$ws.Range("I97").Value = 5761.8
$ws.Range("K97").Value = 5761.8
$ws.Range("M97").Value = -5265.8
$ws.Range("H102").Value = 2627.2
$ws.Range("I102").Value = 2474
$ws.Range("J102").Value = 3240
$ws.Range("K102").Value = 2474
$ws.Range("L102").Value = 3240
$ws.Range("M102").Value = -852
$ws.Range("N102").Value = -6484
$ws.Range("H110").Value = 2425.2307
$ws.Range("I110").Value = 1957.1818
$ws.Range("K110").Value = 1957.1818
$ws.Range("M110").Value = 87.81819999999993

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2111
$ws.Range("I107").Value = 1528.1428
$ws.Range("J107").Value = 3471
$ws.Range("K107").Value = 1528.1428
$ws.Range("L107").Value = 3471
$ws.Range("M107").Value = 391.8571999999999
$ws.Range("N107").Value = -7311

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 78.42856999999999
$ws.Range("I7").Value = 60.666668
$ws.Range("K7").Value = 60.666668
$ws.Range("M7").Value = 52.333332
$ws.Range("H17").Value = 6699.143
$ws.Range("I17").Value = 6815.6665
$ws.Range("J17").Value = 6000
$ws.Range("K17").Value = 6815.6665
$ws.Range("L17").Value = 6000
$ws.Range("M17").Value = -6641.6665
$ws.Range("N17").Value = -6348
$ws.Range("H31").Value = 2523.423
$ws.Range("I31").Value = 2228.3157
$ws.Range("K31").Value = 2228.3157
$ws.Range("M31").Value = -1933.3157
$ws.Range("H34").Value = 2523.423
$ws.Range("I34").Value = 2228.3157
$ws.Range("K34").Value = 2228.3157
$ws.Range("M34").Value = -2026.3157
$ws.Range("H99").Value = 2107.5715
$ws.Range("I99").Value = 1714.5714
$ws.Range("J99").Value = 2500.5715
$ws.Range("K99").Value = 1714.5714
$ws.Range("L99").Value = 2500.5715
$ws.Range("M99").Value = -216.5714
$ws.Range("N99").Value = -5496.5715
$ws.Range("H120").Value = 48364.8
$ws.Range("J120").Value = 48364.8
$ws.Range("L120").Value = 48364.8
$ws.Range("N120").Value = -55622.8
$ws.Range("H126").Value = 2107.5715
$ws.Range("I126").Value = 1714.5714
$ws.Range("J126").Value = 2500.5715
$ws.Range("K126").Value = 5143.7142
$ws.Range("L126").Value = 7501.7145
$ws.Range("M126").Value = -2673.7142
$ws.Range("N126").Value = -12441.7145
$ws.Range("H132").Value = 3367
$ws.Range("I132").Value = 3261.6667
$ws.Range("K132").Value = 9785.000100000001
$ws.Range("M132").Value = -7255.000100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 950.5
$ws.Range("I113").Value = 1046.75
$ws.Range("K113").Value = 3140.25
$ws.Range("M113").Value = -970.25
$ws.Range("H131").Value = 4778107
$ws.Range("I131").Value = 18987.5
$ws.Range("J131").Value = 6681755
$ws.Range("K131").Value = 56962.5
$ws.Range("L131").Value = 20045265
$ws.Range("M131").Value = -51922.5
$ws.Range("N131").Value = -20055345
$ws.Range("H132").Value = 986.3333
$ws.Range("I132").Value = 990
$ws.Range("K132").Value = 8910
$ws.Range("M132").Value = -6380
$ws.Range("H140").Value = 2710
$ws.Range("I140").Value = 1946.4117
$ws.Range("J140").Value = 4873.5
$ws.Range("K140").Value = 5839.2351
$ws.Range("L140").Value = 14620.5
$ws.Range("M140").Value = -659.2350999999999
$ws.Range("N140").Value = -24980.5
$ws.Range("H141").Value = 1503.25
$ws.Range("I141").Value = 1503.25
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 4509.75
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 670.25
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 112
$ws.Range("I2").Value = 115
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 115
$ws.Range("L2").Value = 100
$ws.Range("M2").Value = -2
$ws.Range("N2").Value = -326
$ws.Range("H80").Value = 2651.7812
$ws.Range("I80").Value = 1940.6
$ws.Range("J80").Value = 3279.2942
$ws.Range("K80").Value = 1940.6
$ws.Range("L80").Value = 3279.2942
$ws.Range("M80").Value = -942.5999999999999
$ws.Range("N80").Value = -5275.2942
$ws.Range("H83").Value = 2651.7812
$ws.Range("I83").Value = 1940.6
$ws.Range("J83").Value = 3279.2942
$ws.Range("K83").Value = 9703
$ws.Range("L83").Value = 16396.471
$ws.Range("M83").Value = -4711
$ws.Range("N83").Value = -26380.471
$ws.Range("H113").Value = 236433.22
$ws.Range("I113").Value = 186316.5
$ws.Range("K113").Value = 186316.5
$ws.Range("M113").Value = -184146.5
$ws.Range("H122").Value = 3429.4
$ws.Range("I122").Value = 3074.5
$ws.Range("J122").Value = 3666
$ws.Range("K122").Value = 9223.5
$ws.Range("L122").Value = 10998
$ws.Range("M122").Value = -6773.5
$ws.Range("N122").Value = -15898
$ws.Range("H126").Value = 7765.25
$ws.Range("I126").Value = 7112
$ws.Range("K126").Value = 21336
$ws.Range("M126").Value = -18866

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 1010.5
$ws.Range("I34").Value = 21
$ws.Range("J34").Value = 2000
$ws.Range("K34").Value = 21
$ws.Range("L34").Value = 2000
$ws.Range("M34").Value = 151
$ws.Range("N34").Value = -2344

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2977.6667
$ws.Range("I96").Value = 2685.5715
$ws.Range("K96").Value = 2685.5715
$ws.Range("M96").Value = -1312.5715
$ws.Range("H113").Value = 1177.7273
$ws.Range("I113").Value = 1172.5
$ws.Range("K113").Value = 3517.5
$ws.Range("M113").Value = -1347.5
$ws.Range("H122").Value = 2327.2415
$ws.Range("I122").Value = 2380.1667
$ws.Range("K122").Value = 7140.500100000001
$ws.Range("M122").Value = -4690.500100000001
$ws.Range("H123").Value = 97498.5
$ws.Range("J123").Value = 97498.5
$ws.Range("L123").Value = 97498.5
$ws.Range("N123").Value = -107298.5
$ws.Range("H126").Value = 95176.85000000001
$ws.Range("I126").Value = 133728.5
$ws.Range("K126").Value = 401185.5
$ws.Range("M126").Value = -398715.5
$ws.Range("H136").Value = 4014.6086
$ws.Range("I136").Value = 4375.6313
$ws.Range("J136").Value = 2299.75
$ws.Range("K136").Value = 13126.8939
$ws.Range("L136").Value = 6899.25
$ws.Range("M136").Value = -10576.8939
$ws.Range("N136").Value = -11999.25

Write-Output "Applied all Leve profit updates"